$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Set the value for K10 (mechanism to fetch users' photos task)
$ws.Range("K10").Value = 3

# Reflect the new active cell / selection on the sheet (bottom-right frozen pane)
$ws.Range("A11").Select() | Out-Null
